# Refactor code to save results in a specified folder
# Update the "IPC PO" (column C) values for the data rows, then recompute the
# dependent DELTA (D), DELTA^2 (E) columns and the TOTAL / MSE summary rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newC = @(
    29.81875013623482, 29.96620692625546, 29.64858280488839, 29.67757263402056, 29.16275545037829,
    29.46962369755646, 29.70658430130748, 30.12604549054759, 30.36788804828327, 30.59325418296149,
    30.77393812435264, 31.27779436763707, 31.36702078946621, 31.55713406818036, 31.65189932080112,
    32.01539590513798, 31.84194400932813, 31.89641465392014, 31.85318138283721, 32.48827214235342,
    32.51841802215459, 32.34449417521264, 33.24921042515826, 33.09910619755787, 32.83653587329206,
    33.30289103388547, 33.44744563754524, 33.97102875782325, 34.30255505098672, 34.76595028143137,
    35.25082738242374, 35.47340993938481, 35.78141647081716, 36.31775632249545, 36.86067503638009,
    38.00366613837494, 38.63819033192824, 39.25701150574267, 39.78062723111801, 39.93355248952994,
    39.57900689891961, 40.50373139638501, 40.53828661866802, 41.34804051708238, 40.98852794875395,
    41.12094697642505, 42.2739453400367, 43.48737915769901, 44.38187475288774, 43.93687705961976
)

$firstRow = 2
$lastRow = 51

for ($i = 0; $i -lt $newC.Length; $i++) {
    $row = $firstRow + $i
    $b = $ws.Cells.Item($row, 2).Value2
    $c = $newC[$i]
    $d = $c - $b
    $e = $d * $d

    $ws.Cells.Item($row, 3).Value2 = $c
    $ws.Cells.Item($row, 4).Value2 = $d
    $ws.Cells.Item($row, 5).Value2 = $e
}

# Recompute TOTAL row (52): sum of DELTA and sum of DELTA^2
$totalD = 0
$totalE = 0
for ($row = $firstRow; $row -le $lastRow; $row++) {
    $totalD += $ws.Cells.Item($row, 4).Value2
    $totalE += $ws.Cells.Item($row, 5).Value2
}
$ws.Cells.Item(52, 3).Value2 = $totalD
$ws.Cells.Item(52, 5).Value2 = $totalE

# Recompute MSE row (53): average of DELTA^2
$count = $lastRow - $firstRow + 1
$ws.Cells.Item(53, 5).Value2 = $totalE / $count
